$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row label imported from the perf data source (adds a new shared string)
$ws.Range("A10").Value = "DATA IMPORTED FROM ring.cpp"

# Widen column A so the longer label is readable
$ws.Columns.Item(1).ColumnWidth = 56.33

# Move the active selection to the cell that was just edited
$ws.Range("A10").Select()
